$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "63.888.27"
$ws.Range("E2").Value = "  +0.20%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.621.33"
$ws.Range("E3").Value = "  -0.11%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
Set-TextValue $ws.Range("D5") "594.37"
$ws.Range("E5").Value = "  -0.36%  "

# Row 6
Set-TextValue $ws.Range("D6") "151.26"
$ws.Range("E6").Value = "  +0.84%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.588"
$ws.Range("E8").Value = "  -0.17%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.114"
$ws.Range("E9").Value = "  +4.37%  "

# Row 10
Set-TextValue $ws.Range("D10") "5.81"
$ws.Range("E10").Value = "  +2.01%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.395"
$ws.Range("E11").Value = "  +3.46%  "

# Row 12
$ws.Range("E12").Value = "  +1.09%  "

# Row 13
Set-TextValue $ws.Range("D13") "27.93"
$ws.Range("E13").Value = "  +0.85%  "

# Row 14
Set-TextValue $ws.Range("D14") "3.091.67"
$ws.Range("E14").Value = "  -0.12%  "

# Row 15
Set-TextValue $ws.Range("D15") "63.750.94"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.0000169"
$ws.Range("E16").Value = "  +12.74%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.615.89"
$ws.Range("E17").Value = "  -0.49%  "

# Row 18
Set-TextValue $ws.Range("D18") "12.21"
$ws.Range("E18").Value = "  -0.57%  "

# Row 19
Set-TextValue $ws.Range("D19") "4.78"
$ws.Range("E19").Value = "  +3.07%  "

# Row 20
Set-TextValue $ws.Range("D20") "347.80"
$ws.Range("E20").Value = "  -0.22%  "

# Row 21
Set-TextValue $ws.Range("D21") "7.00"
$ws.Range("E21").Value = "  +1.99%  "

# Row 22
$ws.Range("E22").Value = "  +0.12%  "

# Row 23
Set-TextValue $ws.Range("D23") "67.43"
$ws.Range("E23").Value = "  +1.75%  "

# Row 24
$ws.Range("E24").Value = "  -2.93%  "

# Row 25
$ws.Range("E25").Value = "  +0.57%  "

# Row 26
Set-TextValue $ws.Range("D26") "9.17"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
Set-TextValue $ws.Range("D27") "8.28"
$ws.Range("E27").Value = "  +0.69%  "

# Row 28
Set-TextValue $ws.Range("D28") "548.67"
$ws.Range("E28").Value = "  -2.74%  "

# Row 29
Set-TextValue $ws.Range("D29") "0.162"
$ws.Range("E29").Value = "  -1.82%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.999"
$ws.Range("E30").Value = "  -0.29%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.0₃0908"
$ws.Range("E31").Value = "  +7.67%  "

# Row 32
$ws.Range("E32").Value = "  +0.85%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.82"
$ws.Range("E33").Value = "  +4.52%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.43"
$ws.Range("E34").Value = "  +4.21%  "

# Row 35
Set-TextValue $ws.Range("D35") "6.12"
$ws.Range("E35").Value = "  +0.43%  "

# Row 36
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Range("D36") "0.418"
$ws.Range("E36").Value = "  +2.36%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D37") "164.56"
$ws.Range("E37").Value = "  -2.49%  "

# Row 38
Set-TextValue $ws.Range("D38") "19.96"
$ws.Range("E38").Value = "  +3.21%  "

# Row 39
Set-TextValue $ws.Range("D39") "1.00"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.97"
$ws.Range("E40").Value = "  +1.76%  "

# Row 41
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
Set-TextValue $ws.Range("D42") "167.71"
$ws.Range("E42").Value = "  -1.57%  "

# Row 43
Set-TextValue $ws.Range("D43") "4.09"
$ws.Range("E43").Value = "  +4.62%  "

# Row 44
Set-TextValue $ws.Range("D44") "23.24"
$ws.Range("E44").Value = "  +9.00%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D45") "2.22"
$ws.Range("E45").Value = "  +11.56%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D46") "0.0584"
$ws.Range("E46").Value = "  -1.89%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.637"
$ws.Range("E47").Value = "  +1.27%  "

# Row 48
Set-TextValue $ws.Range("D48") "0.0252"
$ws.Range("E48").Value = "  +1.76%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0971"
$ws.Range("E49").Value = "  +0.27%  "

# Row 50
Set-TextValue $ws.Range("D50") "19.24"
$ws.Range("E50").Value = "  +0.58%  "
